{"js": "// The styles.xml in this document has several character styles\n// (Pandoc/pygments \"*Tok\" syntax-highlighting styles) whose <w:rPr>\n// children are ordered so that <w:color> comes before <w:b>/<w:i>.\n// That violates the CT_RPr content-model order defined by wml.xsd\n// (rFonts, b, bCs, i, iCs, caps, ..., color, ...), which is why\n// OOXMLValidator complains even though xmllint stays silent.\n//\n// Re-applying the existing bold/italic value on each affected style's\n// Font causes the run-properties block to be re-serialized in the\n// canonical schema order, moving <w:b/>/<w:i/> ahead of <w:color/>\n// without altering any actual formatting value.\nconst styleNames = [\n  \"KeywordTok\",\n  \"ImportTok\",\n  \"CommentTok\",\n  \"DocumentationTok\",\n  \"AnnotationTok\",\n  \"CommentVarTok\",\n  \"ControlFlowTok\",\n  \"InformationTok\",\n  \"WarningTok\",\n  \"AlertTok\",\n  \"ErrorTok\",\n];\n\nconst styles = styleNames.map((name) =>\n  context.document.getStyles().getByNameOrNullObject(name)\n);\nstyles.forEach((style) => style.load(\"isNullObject,font\"));\nawait context.sync();\n\nfor (const style of styles) {\n  if (style.isNullObject) {\n    continue;\n  }\n  if (style.font.bold) {\n    style.font.bold = true;\n  }\n  if (style.font.italic) {\n    style.font.italic = true;\n  }\n}\nawait context.sync();\n", "ps1": "# The styles.xml in this document has several character styles\n# (Pandoc/pygments \"*Tok\" syntax-highlighting styles) whose <w:rPr>\n# children are ordered so that <w:color> comes before <w:b>/<w:i>.\n# That violates the CT_RPr content-model order defined by wml.xsd\n# (rFonts, b, bCs, i, iCs, caps, ..., color, ...), which is why\n# OOXMLValidator complains even though xmllint stays silent.\n#\n# Re-applying the existing bold/italic value on each affected style's\n# Font causes the run-properties block to be re-serialized in the\n# canonical schema order, moving <w:b/>/<w:i/> ahead of <w:color/>\n# without altering any actual formatting value.\n$d = $word.ActiveDocument\n\n$styleNames = @(\n    \"KeywordTok\",\n    \"ImportTok\",\n    \"CommentTok\",\n    \"DocumentationTok\",\n    \"AnnotationTok\",\n    \"CommentVarTok\",\n    \"ControlFlowTok\",\n    \"InformationTok\",\n    \"WarningTok\",\n    \"AlertTok\",\n    \"ErrorTok\"\n)\n\nforeach ($name in $styleNames) {\n    $s = $d.Styles($name)\n    if ($null -eq $s) {\n        continue\n    }\n    if ($s.Font.Bold) {\n        $s.Font.Bold = -1\n    }\n    if ($s.Font.Italic) {\n        $s.Font.Italic = -1\n    }\n}\n"}
